$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 533.875
$ws.Range("I33").Value = 547
$ws.Range("K33").Value = 547
$ws.Range("M33").Value = -318
$ws.Range("H34").Value = 1850.1
$ws.Range("I34").Value = 1850.1
$ws.Range("K34").Value = 1850.1
$ws.Range("M34").Value = -1647.1
$ws.Range("H36").Value = 1850.1
$ws.Range("I36").Value = 1850.1
$ws.Range("K36").Value = 1850.1
$ws.Range("M36").Value = -1135.1
$ws.Range("H121").Value = 900
$ws.Range("J121").Value = 900
$ws.Range("L121").Value = 2700
$ws.Range("N121").Value = -6194
$ws.Range("H135").Value = 1129.2667
$ws.Range("I135").Value = 1090.3
$ws.Range("J135").Value = 1207.2
$ws.Range("K135").Value = 9812.699999999999
$ws.Range("L135").Value = 10864.8
$ws.Range("M135").Value = -7277.699999999999
$ws.Range("N135").Value = -15934.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5673.62
$ws.Range("I32").Value = 3927.8696
$ws.Range("K32").Value = 3927.8696
$ws.Range("M32").Value = -3640.8696
$ws.Range("H61").Value = 3886.3142
$ws.Range("I61").Value = 3233.5667
$ws.Range("K61").Value = 3233.5667
$ws.Range("M61").Value = -3021.5667
$ws.Range("H122").Value = 2951.3333
$ws.Range("I122").Value = 2366.3157
$ws.Range("J122").Value = 3961.818
$ws.Range("K122").Value = 7098.9471
$ws.Range("L122").Value = 11885.454
$ws.Range("M122").Value = -4648.9471
$ws.Range("N122").Value = -16785.454
$ws.Range("H132").Value = 3235.973
$ws.Range("I132").Value = 2284.72
$ws.Range("J132").Value = 5217.75
$ws.Range("K132").Value = 6854.16
$ws.Range("L132").Value = 15653.25
$ws.Range("M132").Value = -4324.16
$ws.Range("N132").Value = -20713.25
$ws.Range("H136").Value = 3886.3142
$ws.Range("I136").Value = 3233.5667
$ws.Range("K136").Value = 9700.7001
$ws.Range("M136").Value = -7150.7001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4790.1113
$ws.Range("I86").Value = 4434.6665
$ws.Range("K86").Value = 4434.6665
$ws.Range("M86").Value = -3311.6665
$ws.Range("H89").Value = 4790.1113
$ws.Range("I89").Value = 4434.6665
$ws.Range("K89").Value = 22173.3325
$ws.Range("M89").Value = -16557.3325
$ws.Range("H105").Value = 10047.128
$ws.Range("I105").Value = 8164.9355
$ws.Range("J105").Value = 17340.625
$ws.Range("K105").Value = 8164.9355
$ws.Range("L105").Value = 17340.625
$ws.Range("M105").Value = -6417.9355
$ws.Range("N105").Value = -20834.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1442.08
$ws.Range("I16").Value = 732.7143
$ws.Range("K16").Value = 732.7143
$ws.Range("M16").Value = -445.7143
$ws.Range("H22").Value = 4562.8
$ws.Range("I22").Value = 988
$ws.Range("J22").Value = 9925
$ws.Range("K22").Value = 988
$ws.Range("L22").Value = 9925
$ws.Range("M22").Value = -638
$ws.Range("N22").Value = -10625
$ws.Range("H37").Value = 33333.332
$ws.Range("J37").Value = 34000
$ws.Range("L37").Value = 34000
$ws.Range("N37").Value = -34214
$ws.Range("H113").Value = 1442.08
$ws.Range("I113").Value = 732.7143
$ws.Range("K113").Value = 732.7143
$ws.Range("M113").Value = 1437.2857
$ws.Range("H132").Value = 3266.634
$ws.Range("I132").Value = 2765.1614
$ws.Range("J132").Value = 4821.2
$ws.Range("K132").Value = 8295.484199999999
$ws.Range("L132").Value = 14463.6
$ws.Range("M132").Value = -5765.484199999999
$ws.Range("N132").Value = -19523.6
$ws.Range("H134").Value = 2431.5
$ws.Range("I134").Value = 1558.8
$ws.Range("J134").Value = 4855.6665
$ws.Range("K134").Value = 4676.4
$ws.Range("L134").Value = 14566.9995
$ws.Range("M134").Value = -2141.4
$ws.Range("N134").Value = -19636.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4999.9
$ws.Range("I55").Value = 1715.9
$ws.Range("J55").Value = 8283.9
$ws.Range("K55").Value = 5147.700000000001
$ws.Range("L55").Value = 24851.7
$ws.Range("M55").Value = -4970.700000000001
$ws.Range("N55").Value = -25205.7
$ws.Range("H129").Value = 15158697
$ws.Range("I129").Value = 1594.8
$ws.Range("J129").Value = 27789616
$ws.Range("K129").Value = 4784.4
$ws.Range("L129").Value = 83368848
$ws.Range("M129").Value = 215.6000000000004
$ws.Range("N129").Value = -83378848
$ws.Range("H131").Value = 6798438.5
$ws.Range("J131").Value = 4987650.5
$ws.Range("L131").Value = 14962951.5
$ws.Range("N131").Value = -14973031.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 29333.334
$ws.Range("J96").Value = 29333.334
$ws.Range("L96").Value = 29333.334
$ws.Range("N96").Value = -34825.334
$ws.Range("H138").Value = 76291.25
$ws.Range("J138").Value = 76291.25
$ws.Range("L138").Value = 76291.25
$ws.Range("N138").Value = -86571.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9097057
$ws.Range("I7").Value = 11768927
$ws.Range("K7").Value = 11768927
$ws.Range("M7").Value = -11768815
$ws.Range("H16").Value = 1755
$ws.Range("I16").Value = 577.2857
$ws.Range("J16").Value = 9999
$ws.Range("K16").Value = 577.2857
$ws.Range("L16").Value = 9999
$ws.Range("M16").Value = -407.2857
$ws.Range("N16").Value = -10339
$ws.Range("H40").Value = 2807105.2
$ws.Range("I40").Value = 3667881.2
$ws.Range("K40").Value = 3667881.2
$ws.Range("M40").Value = -3667745.2
$ws.Range("H55").Value = 1854209
$ws.Range("J55").Value = 5190.5
$ws.Range("L55").Value = 5190.5
$ws.Range("N55").Value = -5536.5
$ws.Range("H61").Value = 8903.916999999999
$ws.Range("I61").Value = 7760.3335
$ws.Range("J61").Value = 12334.667
$ws.Range("K61").Value = 7760.3335
$ws.Range("L61").Value = 12334.667
$ws.Range("M61").Value = -7558.3335
$ws.Range("N61").Value = -12738.667
$ws.Range("H93").Value = 1515.0714
$ws.Range("I93").Value = 1704.5238
$ws.Range("J93").Value = 946.7143
$ws.Range("K93").Value = 1704.5238
$ws.Range("L93").Value = 946.7143
$ws.Range("M93").Value = -456.5237999999999
$ws.Range("N93").Value = -3442.7143
$ws.Range("H113").Value = 8903.916999999999
$ws.Range("I113").Value = 7760.3335
$ws.Range("J113").Value = 12334.667
$ws.Range("K113").Value = 7760.3335
$ws.Range("L113").Value = 12334.667
$ws.Range("M113").Value = -5590.3335
$ws.Range("N113").Value = -16674.667
$ws.Range("H126").Value = 9097057
$ws.Range("I126").Value = 11768927
$ws.Range("K126").Value = 35306781
$ws.Range("M126").Value = -35304311

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1686.0227
$ws.Range("I122").Value = 1216.081
$ws.Range("K122").Value = 3648.242999999999
$ws.Range("M122").Value = -1198.242999999999
$ws.Range("H132").Value = 2725.484
$ws.Range("J132").Value = 5232.8335
$ws.Range("L132").Value = 15698.5005
$ws.Range("N132").Value = -20758.5005
$ws.Range("H136").Value = 2539.9048
$ws.Range("I136").Value = 877.9259
$ws.Range("J136").Value = 5531.467
$ws.Range("K136").Value = 2633.7777
$ws.Range("L136").Value = 16594.401
$ws.Range("M136").Value = -83.77769999999964
$ws.Range("N136").Value = -21694.401
